# Database schema migration: append one new reading row to each sheet.

$wb = $excel.ActiveWorkbook

function Add-Row {
    param(
        $ws,
        [int]$row,
        [double]$a,
        [string]$b,
        [string]$c,
        [string]$d,
        [string]$e,
        [double]$f,
        $g,
        [double]$h,
        [double]$i
    )

    # Column A: timestamp, same date/time style as the row above it.
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 1).Value = $a

    # Columns B-E: hex-ish payload strings, stored as text.
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e

    # Column F: plain numeric.
    $ws.Cells.Item($row, 6).Value = $f

    # Column G: numeric, unless it must be preserved verbatim as text
    # (huge integer literal that loses precision as a double).
    if ($g -is [string]) {
        $ws.Cells.Item($row, 7).NumberFormat = "@"
        $ws.Cells.Item($row, 7).Value = $g
    } else {
        $ws.Cells.Item($row, 7).Value = $g
    }

    # Columns H-I: plain numeric.
    $ws.Cells.Item($row, 8).Value = $h
    $ws.Cells.Item($row, 9).Value = $i
}

# PowerShell parser here chokes on bare scientific-notation literals
# (5.68631262647114e+23) -- cast the string form to [double] instead.
$bigNum = [double]"5.68631262647114e+23"

# Sheet 1: ROW50-FE-LIFTER -- add row 90
$ws1 = $wb.Worksheets.Item("ROW50-FE-LIFTER")
Add-Row $ws1 90 45769.29095546297 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x3e" "0xe" 400 $bigNum 318 14

# Sheet 2: ROW50-MID-LIFTER -- add row 92
$ws2 = $wb.Worksheets.Item("ROW50-MID-LIFTER")
Add-Row $ws2 92 45769.25613425926 "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x46" "0x19" 400 "568631262647113771663628" 326 25

# Sheet 3: ROW11-FE-LIFTER -- add row 90
$ws3 = $wb.Worksheets.Item("ROW11-FE-LIFTER")
Add-Row $ws3 90 45769.32291275463 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x01,0x3e" "0x14" 400 $bigNum 318 20

# Sheet 4: ROW11-MID-LIFTER -- add row 90
$ws4 = $wb.Worksheets.Item("ROW11-MID-LIFTER")
Add-Row $ws4 90 45769.44618920139 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," "0x01,0x46" "0x19" 400 $bigNum 326 25

Write-Output "migration applied"
